# feat: add 2022-Q1 data
#
# - Insert a new sheet "2022-Q1" (fund-holding detail) positioned between
#   "2021-Q4" and "总计", built as a copy of "2021-Q4" so it inherits the
#   same header/row styling, then re-populated with the new quarter's data.
# - Update the "总计" (totals) sheet: prepend a new summary row for
#   "2022-Q1" and shift the existing "2021-Q4" / "2021-Q3" rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q1" detail sheet by copying "2021-Q4" (keeps the
#    same column widths / header style / border style as the sibling
#    quarter sheets) and placing the copy right after it.
# ---------------------------------------------------------------------
$srcQ = $wb.Worksheets.Item("2021-Q4")
$srcQ.Copy($null, $srcQ)
$q1 = $wb.Worksheets.Item($srcQ.Index + 1)
$q1.Name = "2022-Q1"

# Column D header changes from "基金规模" (inherited from 2021-Q4) -> stays
# "基金规模" per the diff, so no header text changes are required; only the
# data rows need to be rewritten/expanded to 6 funds.

# Extend the style used by column A (row index marker, e.g. s="2") down to
# the additional rows 4-7 (the source sheet only had rows 2-3).
$q1.Range("A2").Copy()
$q1.Range("A4:A7").PasteSpecial(-4122)

# Force columns B:G to be stored as text (so values like "006529" keep
# their leading zero and decimals like "43.69" aren't coerced to numbers),
# then drop back to the default "Normal" style so no stray formatting is
# left behind.
$q1.Range("B2:G7").NumberFormat = "@"

$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "006529"
$q1.Range("C2").Value = "中欧匠心两年持有期混合A"
$q1.Range("D2").Value = "43.69"
$q1.Range("E2").Value = "88.30"
$q1.Range("F2").Value = "2.35"
$q1.Range("G2").Value = "1.0267"
$q1.Range("H2").Value = 10

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "006530"
$q1.Range("C3").Value = "中欧匠心两年持有期混合C"
$q1.Range("D3").Value = "4.59"
$q1.Range("E3").Value = "88.30"
$q1.Range("F3").Value = "2.35"
$q1.Range("G3").Value = "0.1079"
$q1.Range("H3").Value = 10

$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "011651"
$q1.Range("C4").Value = "招商港股通核心精选股票A"
$q1.Range("D4").Value = "2.81"
$q1.Range("E4").Value = "81.27"
$q1.Range("F4").Value = "2.34"
$q1.Range("G4").Value = "0.0658"
$q1.Range("H4").Value = 10

$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "004266"
$q1.Range("C5").Value = "招商沪港深科技创新主题精选灵活配置混合A"
$q1.Range("D5").Value = "1.29"
$q1.Range("E5").Value = "88.85"
$q1.Range("F5").Value = "2.01"
$q1.Range("G5").Value = "0.0259"
$q1.Range("H5").Value = 10

$q1.Range("A6").Value = 4
$q1.Range("B6").Value = "011652"
$q1.Range("C6").Value = "招商港股通核心精选股票C"
$q1.Range("D6").Value = "0.94"
$q1.Range("E6").Value = "81.27"
$q1.Range("F6").Value = "2.34"
$q1.Range("G6").Value = "0.0220"
$q1.Range("H6").Value = 10

$q1.Range("A7").Value = 5
$q1.Range("B7").Value = "010754"
$q1.Range("C7").Value = "招商沪港深科技创新主题精选灵活配置混合C"
$q1.Range("D7").Value = "0.28"
$q1.Range("E7").Value = "88.85"
$q1.Range("F7").Value = "2.01"
$q1.Range("G7").Value = "0.0056"
$q1.Range("H7").Value = 10

$q1.Range("B2:G7").Style = "Normal"

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: add a new top data row for 2022-Q1 and push
#    the prior two quarters down by one row.
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

$oldB2 = $zj.Range("B2").Value()
$oldC2 = $zj.Range("C2").Value()
$oldD2 = $zj.Range("D2").Value()
$oldB3 = $zj.Range("B3").Value()
$oldC3 = $zj.Range("C3").Value()
$oldD3 = $zj.Range("D3").Value()

# Row 4 is brand new -- copy row 2's column-A style down to it first.
$zj.Range("A2").Copy()
$zj.Range("A4").PasteSpecial(-4122)

$zj.Range("A4").Value = 2
$zj.Range("B4").Value = $oldB3
$zj.Range("C4").Value = $oldC3
$zj.Range("D4").Value = $oldD3

$zj.Range("A3").Value = 1
$zj.Range("B3").Value = $oldB2
$zj.Range("C3").Value = $oldC2
$zj.Range("D3").Value = $oldD2

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 6
$zj.Range("D2").Value = 1.25
